$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (existing) values to the new strings
$ws.Range("A2").Value = " Exam_Prerequisit_for_iProc_TC_ID_79"
$ws.Range("B2").Value = "@iProctorRegression Verify Elumina Create Exam with survey section"
$ws.Range("C2").Value = "passed"

# Add new row 3
$ws.Range("A3").Value = " Exam_Prerequisit_for_iProc_TC_ID_79"
$ws.Range("B3").Value = "@iProctorRegression Verify Elumina Registration"
$ws.Range("C3").Value = "passed"

# Add new row 4
$ws.Range("A4").Value = "iProc_TC_ID_79"
$ws.Range("B4").Value = "@iProctorRegression Verify Validation of Candidate answering survey questions"
$ws.Range("C4").Value = "passed"
